$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.188.33'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = '1.581.77'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = "'209.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D6').Value = "'0.495"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('D8').Value = "'0.0610"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('D9').Value = "'0.245"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('D10').Value = "'19.49"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('D11').Value = "'0.0846"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '1.804.85'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'4.05"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.575.21'
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('D15').Value = "'0.515"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = "'64.46"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = '26.202.18'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('D18').Value = '0.0₃0734'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = "'7.27"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').Value = "'207.11"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').Value = "'4.25"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D24').Value = "'8.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('D25').Value = "'144.94"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').Value = "'7.01"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').Value = "'15.20"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('D32').Value = "'3.21"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.38%  '
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').Value = '1.279.15'
$ws.Range('E34').Value = '  -1.10%  '
$ws.Range('D35').Value = "'2.46"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'0.610"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.94%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = "'1.21"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.35%  '
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').Value = "'0.813"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('D41').Value = "'5.57"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.35%  '
$ws.Range('D42').Value = "'0.767"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('D44').Value = "'62.24"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').Value = '1.718.50'
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('D46').Value = "'88.84"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.91%  '
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('E49').Value = '  -1.74%  '
$ws.Range('E50').Value = '  -0.26%  '
$ws.Range('D51').Value = '0.0₇0955'
$ws.Range('E51').Value = '  -10.19%  '
